$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PersonalizerItems")

# Add the new row 20, mirroring row 19 but with the "B" variant selector values
$ws.Range("A20").Value = "https://www.microsoft.com/en-us/store/collections/certifiedrefurbishedsurface"
$ws.Range("B20").Value = "28806B"
$ws.Range("C20").Value = ".wrapper-28806B .col"
$ws.Range("D20").Value = "0.6,0.3,0.1"
$ws.Range("E20").Value = $false

# Default row height changes from 15.5 to 14.5, while all existing data rows
# keep an explicit 15.5 row height (so they render unchanged)
$ws.Rows.Item(1).RowHeight = 15.5
for ($r = 1; $r -le 20; $r++) {
    $ws.Rows.Item($r).RowHeight = 15.5
}
$ws.StandardHeight = 14.5

# Update ignored errors region to cover the full new range as a single block
$ws.Range("A1:E20").ErrorCheckingOptions.NumberAsText = $true

# Reset the view: scroll back to the top, and select B1
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B1").Select()
